# Excel COM-interop script applying a market-price / leve-profit data refresh
# across the per-job Leve tables (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Values come from a scheduled data-refresh run (currentAveragePrice* / LevePrice*
# / LeveProfit* columns); a few rows also gain/lose a trailing HQ-profit cell
# depending on whether that SKU has an HQ variant priced this run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 310
$ws.Range("I58").Value = 310
$ws.Range("K58").Value = 930
$ws.Range("M58").Value = -780

$ws.Range("H64").Value = 15416.5
$ws.Range("I64").Value = 4166.3335
$ws.Range("K64").Value = 4166.3335
$ws.Range("M64").Value = -3918.3335

$ws.Range("H67").Value = 15416.5
$ws.Range("I67").Value = 4166.3335
$ws.Range("K67").Value = 4166.3335
$ws.Range("M67").Value = -3308.3335

$ws.Range("H69").Value = 7974.875
$ws.Range("I69").Value = 2000
$ws.Range("K69").Value = 6000
$ws.Range("M69").Value = -5126

$ws.Range("H72").Value = 7974.875
$ws.Range("I72").Value = 2000
$ws.Range("K72").Value = 18000
$ws.Range("M72").Value = -13632

$ws.Range("H86").Value = 5302.25
$ws.Range("I86").Value = 4237.0835
$ws.Range("J86").Value = 6900
$ws.Range("K86").Value = 4237.0835
$ws.Range("L86").Value = 6900
$ws.Range("M86").Value = -3114.0835
$ws.Range("N86").Value = -9146

$ws.Range("H89").Value = 5302.25
$ws.Range("I89").Value = 4237.0835
$ws.Range("J89").Value = 6900
$ws.Range("K89").Value = 21185.4175
$ws.Range("L89").Value = 34500
$ws.Range("M89").Value = -15569.4175
$ws.Range("N89").Value = -45732

$ws.Range("H116").Value = 15240.667
$ws.Range("I116").Value = 6298.875
$ws.Range("J116").Value = 33124.25
$ws.Range("K116").Value = 6298.875
$ws.Range("L116").Value = 33124.25
$ws.Range("M116").Value = -2856.875
$ws.Range("N116").Value = -40008.25

$ws.Range("H136").Value = 60496
$ws.Range("J136").Value = 60496
$ws.Range("L136").Value = 60496
$ws.Range("N136").Value = -70696

$ws.Range("H137").Value = 29363.523
$ws.Range("I137").Value = 18010.756
$ws.Range("K137").Value = 54032.268
$ws.Range("M137").Value = -51482.268

$ws.Range("H138").Value = 4739.4653
$ws.Range("I138").Value = 2783.1667
$ws.Range("J138").Value = 4886.1875
$ws.Range("K138").Value = 8349.500100000001
$ws.Range("L138").Value = 14658.5625
$ws.Range("M138").Value = -3209.500100000001
$ws.Range("N138").Value = -24938.5625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 5607.5
$ws.Range("I21").Value = 5607.5
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 5607.5
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -5233.5
$ws.Range("N21").ClearContents()

$ws.Range("H32").Value = 2666.4219
$ws.Range("I32").Value = 1666.9482
$ws.Range("J32").Value = 12328
$ws.Range("K32").Value = 1666.9482
$ws.Range("L32").Value = 12328
$ws.Range("M32").Value = -1379.9482
$ws.Range("N32").Value = -12902

$ws.Range("H61").Value = 2054.068
$ws.Range("I61").Value = 1909.125
$ws.Range("J61").Value = 3503.5
$ws.Range("K61").Value = 1909.125
$ws.Range("L61").Value = 3503.5
$ws.Range("M61").Value = -1697.125
$ws.Range("N61").Value = -3927.5

$ws.Range("H110").Value = 12427.096
$ws.Range("I110").Value = 16580.75
$ws.Range("K110").Value = 16580.75
$ws.Range("M110").Value = -14535.75

$ws.Range("H122").Value = 3943.4482
$ws.Range("I122").Value = 3671.25
$ws.Range("K122").Value = 11013.75
$ws.Range("M122").Value = -8563.75

$ws.Range("H125").Value = 40715
$ws.Range("J125").Value = 40715
$ws.Range("L125").Value = 40715
$ws.Range("N125").Value = -50555

$ws.Range("H132").Value = 4243.579
$ws.Range("I132").Value = 3656.516
$ws.Range("K132").Value = 10969.548
$ws.Range("M132").Value = -8439.548000000001

$ws.Range("H136").Value = 2054.068
$ws.Range("I136").Value = 1909.125
$ws.Range("J136").Value = 3503.5
$ws.Range("K136").Value = 5727.375
$ws.Range("L136").Value = 10510.5
$ws.Range("M136").Value = -3177.375
$ws.Range("N136").Value = -15610.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2284.6667
$ws.Range("I107").Value = 1873.8948
$ws.Range("K107").Value = 1873.8948
$ws.Range("M107").Value = 46.10519999999997

$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 44486.633
$ws.Range("I31").Value = 53142.562
$ws.Range("J31").Value = 13710
$ws.Range("K31").Value = 53142.562
$ws.Range("L31").Value = 13710
$ws.Range("M31").Value = -52847.562
$ws.Range("N31").Value = -14300

$ws.Range("H34").Value = 44486.633
$ws.Range("I34").Value = 53142.562
$ws.Range("J34").Value = 13710
$ws.Range("K34").Value = 53142.562
$ws.Range("L34").Value = 13710
$ws.Range("M34").Value = -52940.562
$ws.Range("N34").Value = -14114

$ws.Range("H41").Value = 11500
$ws.Range("J41").Value = 11000
$ws.Range("L41").Value = 11000
$ws.Range("N41").Value = -11856

$ws.Range("H50").Value = 6933.3335

$ws.Range("H58").Value = 2286.926
$ws.Range("I58").Value = 2075.2104
$ws.Range("K58").Value = 2075.2104
$ws.Range("M58").Value = -1872.2104

$ws.Range("H59").Value = 25600
$ws.Range("J59").Value = 25600
$ws.Range("L59").Value = 25600
$ws.Range("N59").Value = -27890

$ws.Range("H60").Value = 14030.333

$ws.Range("H122").Value = 1666.5
$ws.Range("I122").Value = 1666.5
$ws.Range("K122").Value = 4999.5
$ws.Range("M122").Value = -2549.5

$ws.Range("H136").Value = 2286.926
$ws.Range("I136").Value = 2075.2104
$ws.Range("K136").Value = 6225.6312
$ws.Range("M136").Value = -3675.6312

$ws.Range("H141").Value = 66108.57000000001
$ws.Range("J141").Value = 66108.57000000001
$ws.Range("L141").Value = 66108.57000000001
$ws.Range("N141").Value = -76468.57000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 350
$ws.Range("I22").Value = 475.5
$ws.Range("K22").Value = 1426.5
$ws.Range("M22").Value = -1257.5

$ws.Range("H23").Value = 137.54546
$ws.Range("I23").Value = 134.77777
$ws.Range("J23").Value = 150
$ws.Range("K23").Value = 404.33331
$ws.Range("L23").Value = 450
$ws.Range("M23").Value = -169.33331
$ws.Range("N23").Value = -920

$ws.Range("H27").Value = 350
$ws.Range("I27").Value = 475.5
$ws.Range("K27").Value = 1426.5
$ws.Range("M27").Value = -1324.5

$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()

$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

$ws.Range("H58").Value = 5904.25
$ws.Range("I58").Value = 5329
$ws.Range("J58").Value = 6479.5
$ws.Range("K58").Value = 15987
$ws.Range("L58").Value = 19438.5
$ws.Range("M58").Value = -15859
$ws.Range("N58").Value = -19694.5

$ws.Range("H107").Value = 1330.7693
$ws.Range("J107").Value = 2666.6
$ws.Range("L107").Value = 7999.799999999999
$ws.Range("N107").Value = -11839.8

$ws.Range("H113").Value = 294.7
$ws.Range("I113").Value = 336.4737
$ws.Range("J113").Value = 222.54546
$ws.Range("K113").Value = 1009.4211
$ws.Range("L113").Value = 667.6363799999999
$ws.Range("M113").Value = 1160.5789
$ws.Range("N113").Value = -5007.63638

$ws.Range("H124").Value = 20103.656
$ws.Range("J124").Value = 20103.656
$ws.Range("L124").Value = 60310.96799999999
$ws.Range("N124").Value = -70130.96799999999

$ws.Range("H129").Value = 4501962.5
$ws.Range("I129").Value = 9900782
$ws.Range("J129").Value = 2945.8333
$ws.Range("K129").Value = 29702346
$ws.Range("L129").Value = 8837.499899999999
$ws.Range("M129").Value = -29697346
$ws.Range("N129").Value = -18837.4999

$ws.Range("H131").Value = 18009.363
$ws.Range("I131").Value = 101132.2
$ws.Range("J131").Value = 3166
$ws.Range("K131").Value = 303396.6
$ws.Range("L131").Value = 9498
$ws.Range("M131").Value = -298356.6
$ws.Range("N131").Value = -19578

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 17284.771
$ws.Range("I132").Value = 18505.219
$ws.Range("K132").Value = 55515.65700000001
$ws.Range("M132").Value = -52985.65700000001

$ws.Range("H135").Value = 59999
$ws.Range("J135").Value = 59999
$ws.Range("L135").Value = 59999
$ws.Range("N135").Value = -70139

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H94").Value = 50330
$ws.Range("J94").Value = 50330
$ws.Range("L94").Value = 50330
$ws.Range("N94").Value = -51682

$ws.Range("H122").Value = 405709.72
$ws.Range("I122").Value = 592890.4
$ws.Range("K122").Value = 1778671.2
$ws.Range("M122").Value = -1776221.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 15004998

$ws.Range("H17").Value = 605
$ws.Range("I17").Value = 10
$ws.Range("J17").Value = 1200
$ws.Range("K17").Value = 10
$ws.Range("L17").Value = 1200
$ws.Range("M17").Value = 162
$ws.Range("N17").Value = -1544

$ws.Range("H22").Value = 7816.6665
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 9180
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 9180
$ws.Range("M22").Value = -707
$ws.Range("N22").Value = -9766

$ws.Range("H81").Value = 7193.76
$ws.Range("I81").Value = 9720.385
$ws.Range("J81").Value = 4456.5835
$ws.Range("K81").Value = 19440.77
$ws.Range("L81").Value = 8913.166999999999
$ws.Range("M81").Value = -18379.77
$ws.Range("N81").Value = -11035.167

$ws.Range("H84").Value = 7193.76
$ws.Range("I84").Value = 9720.385
$ws.Range("J84").Value = 4456.5835
$ws.Range("K84").Value = 97203.85000000001
$ws.Range("L84").Value = 44565.835
$ws.Range("M84").Value = -91899.85000000001
$ws.Range("N84").Value = -55173.835

$ws.Range("H136").Value = 3193.5
$ws.Range("I136").Value = 2891.6428
$ws.Range("K136").Value = 8674.928400000001
$ws.Range("M136").Value = -6124.928400000001
